$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solutions")

# Insert a new column before column I (MST results shift from I -> J,
# %Over formulas shift from J -> K; Excel auto-fixes formula references
# and the conditional-formatting / dimension ranges when inserting).
$ws.Columns.Item(9).Insert()

# Header for the newly inserted column.
$ws.Range("I1").Value = "Time"

# New per-instance MST running times (seconds) for the inserted column.
$ws.Range("I2").Value = 0.01
$ws.Range("I3").Value = 0.224
$ws.Range("I4").Value = 0.091
$ws.Range("I5").Value = 0.249
$ws.Range("I6").Value = 0.002
$ws.Range("I7").Value = 0.795
$ws.Range("I8").Value = 0.627
$ws.Range("I9").Value = 0.045
$ws.Range("I10").Value = 42.779
$ws.Range("I11").Value = 1.492
$ws.Range("I12").Value = 3.147
$ws.Range("I13").Value = 0.002
$ws.Range("I14").Value = 1.846

# Updated MST solution values (now in column J) for the rows whose
# results changed with this run.
$ws.Range("J10").Value = 808235
$ws.Range("J11").Value = 1060717
$ws.Range("J12").Value = 1618300
$ws.Range("J13").Value = 65561
